{"js": "// Remove the trailing \"Ver no Jupiter...\" / \"\u00a9 2020 ...\" footer block\n// (and the blank paragraph right before it) that used to follow the\n// last bibliography entry (\"SERAFINI, Maria Jos\u00e9...\").\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targetTexts = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\",\n];\n\nconst items = paragraphs.items;\n\n// Find the \"Ver no Jupiter...\" paragraph; the blank paragraph that\n// immediately precedes it (left over from the footer block) and the\n// copyright paragraph that immediately follows it are removed together.\nlet jupiterIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === targetTexts[0]) {\n    jupiterIndex = i;\n    break;\n  }\n}\n\nconst toDelete = [];\nif (jupiterIndex !== -1) {\n  if (jupiterIndex - 1 >= 0 && items[jupiterIndex - 1].text === \"\") {\n    toDelete.push(items[jupiterIndex - 1]);\n  }\n  toDelete.push(items[jupiterIndex]);\n  if (\n    jupiterIndex + 1 < items.length &&\n    items[jupiterIndex + 1].text === targetTexts[1]\n  ) {\n    toDelete.push(items[jupiterIndex + 1]);\n  }\n}\n\nfor (const p of toDelete) {\n  p.delete();\n}\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" / \"\u00a9 2020 ...\" footer block\n# (and the blank paragraph right before it) that used to follow the\n# last bibliography entry (\"SERAFINI, Maria Jos\u00e9...\").\n$d = $word.ActiveDocument\n\n$jupiterText = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$copyrightText = \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n\n$count = $d.Paragraphs.Count\n$jupiterIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $text = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)\n    if ($text -eq $jupiterText) {\n        $jupiterIndex = $i\n        break\n    }\n}\n\nif ($jupiterIndex -ge 1) {\n    $indexesToDelete = New-Object System.Collections.ArrayList\n\n    if ($jupiterIndex + 1 -le $count) {\n        $nextText = $d.Paragraphs.Item($jupiterIndex + 1).Range.Text.TrimEnd([char]13, [char]7)\n        if ($nextText -eq $copyrightText) {\n            [void]$indexesToDelete.Add($jupiterIndex + 1)\n        }\n    }\n\n    [void]$indexesToDelete.Add($jupiterIndex)\n\n    if ($jupiterIndex - 1 -ge 1) {\n        $prevText = $d.Paragraphs.Item($jupiterIndex - 1).Range.Text.TrimEnd([char]13, [char]7)\n        if ($prevText -eq \"\") {\n            [void]$indexesToDelete.Add($jupiterIndex - 1)\n        }\n    }\n\n    # Delete from the highest index down so earlier indices stay valid.\n    $sorted = $indexesToDelete | Sort-Object -Descending\n    foreach ($idx in $sorted) {\n        $d.Paragraphs.Item($idx).Range.Delete()\n    }\n}\n"}
